$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Rep ID"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 11
